# Updated the stats after last release
# Adds one new day's worth of stats (row 81) to the "Data" table, mirroring
# the same formulas / number-formats used by the previous last row (80),
# and tidies up the formatting artifact Excel leaves behind on the row that
# used to be last.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$lo = $ws.ListObjects.Item("Data")

# --- 1. Grow the table by one row -----------------------------------------
$newListRow = $lo.ListRows.Add()
$r = $newListRow.Range.Row   # should be 81

# --- 2. Carry the formatting of the previous last row into the new row ----
$ws.Range("A80:AK80").Copy()
$ws.Range("A81:AK81").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Fill in the new row's data -----------------------------------------
$ws.Cells.Item($r, 1).Value2 = 45230     # Date
$ws.Cells.Item($r, 2).Value2 = 345       # Stars
$ws.Cells.Item($r, 3).Value2 = 129       # Forks
$ws.Cells.Item($r, 4).Value2 = 113       # Authors
$ws.Cells.Item($r, 5).Value2 = 281       # Versions
$ws.Cells.Item($r, 6).Value2 = 232       # GH Releases
$ws.Cells.Item($r, 7).Value2 = 5728      # LoC
$ws.Cells.Item($r, 8).Formula = "=Data[[#This Row],[LoC]]-G80"              # ∆LoC
$ws.Cells.Item($r, 9).Value2 = 6867      # Shell
$ws.Cells.Item($r, 10).Value2 = 1967     # MD
$ws.Cells.Item($r, 11).Value2 = 567      # YAML
$ws.Cells.Item($r, 12).Value2 = 290      # Text
$ws.Cells.Item($r, 13).Value2 = 134      # make
$ws.Cells.Item($r, 14).Value2 = 60       # Bash
$ws.Cells.Item($r, 15).Value2 = 16       # ini
$ws.Cells.Item($r, 16).Formula = "=SUM(Data[[#This Row],[Shell]:[Bash]])"   # Total
$ws.Cells.Item($r, 17).Formula = "=Data[[#This Row],[Total]]-P80"           # ∆Total
$ws.Cells.Item($r, 18).Value2 = 2136     # Commits
$ws.Cells.Item($r, 19).Value2 = 4536     # File Changes
$ws.Cells.Item($r, 20).Value2 = 71120    # Insertions
$ws.Cells.Item($r, 21).Value2 = 48752    # Deletions
$ws.Cells.Item($r, 22).Value2 = 2        # Open issues
$ws.Cells.Item($r, 23).Value2 = 1        # Open bugs
$ws.Cells.Item($r, 24).Value2 = 274      # Closed issues
$ws.Cells.Item($r, 25).Formula = "=Data[[#This Row],[Open issues]]+Data[[#This Row],[Closed issues]]"   # Issues
$ws.Cells.Item($r, 26).Value2 = 0        # Open pull requests
$ws.Cells.Item($r, 27).Value2 = 176      # Closed pull requests
$ws.Cells.Item($r, 28).Formula = "=Data[[#This Row],[Open pull requests]]+Data[[#This Row],[Closed pull requests]]"   # Pull requests
$ws.Cells.Item($r, 29).Value2 = 159      # Command line options
$ws.Cells.Item($r, 30).Value2 = 165      # Tests
$ws.Cells.Item($r, 31).Value2 = 7        # GH workflows
$ws.Cells.Item($r, 32).Value2 = 0        # Running
$ws.Cells.Item($r, 33).Value2 = 127      # Failed
$ws.Cells.Item($r, 34).Value2 = 994      # OK
$ws.Cells.Item($r, 35).Value2 = 7        # Cancelled
                                          # GH runs (col 36) intentionally left blank
$ws.Cells.Item($r, 37).Formula = "=SUM(Data[[#This Row],[Running]:[GH runs]])"   # ∑runs

# --- 4. The row that used to be last loses the stray "General" style that
#        Excel had tagged it with while it was the table's edge row --------
$ws.Range("B79:C79").Copy()
$ws.Range("B80:C80").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 5. Leave the selection where Excel would after typing the new row ----
$ws.Range("AL81").Select()
